$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "29.388.94"
Set-TextValue "E2" "  -0.04%  "
Set-TextValue "D3" "1.849.45"
Set-TextValue "E3" "  +0.03%  "
Set-TextValue "D4" "0.9992"
Set-TextValue "E4" "  -0.04%  "
Set-TextValue "D5" "240.25"
Set-TextValue "D6" "0.6283"
Set-TextValue "E6" "  -0.25%  "
Set-TextValue "D7" "1.000"
Set-TextValue "E7" "  -0.02%  "
Set-TextValue "E8" "  +0.05%  "
Set-TextValue "E9" "  -1.13%  "
Set-TextValue "D10" "24.72"
Set-TextValue "E10" "  +0.99%  "
Set-TextValue "E11" "  -0.04%  "
Set-TextValue "D12" "5.031"
Set-TextValue "E12" "  +0.42%  "
Set-TextValue "D13" "0.6782"
Set-TextValue "E14" "  -2.52%  "
Set-TextValue "D15" "83.27"
Set-TextValue "E15" "  -0.31%  "
Set-TextValue "D16" "6.158"
Set-TextValue "E16" "  +0.21%  "
Set-TextValue "D17" "29.432.58"
Set-TextValue "E17" "  +0.00%  "
Set-TextValue "D18" "227.68"
Set-TextValue "E18" "  -0.41%  "
Set-TextValue "D20" "0.9996"
Set-TextValue "E20" "  -0.07%  "
Set-TextValue "D21" "7.512"
Set-TextValue "E21" "  +0.78%  "
Set-TextValue "E22" "  -0.04%  "
Set-TextValue "D23" "158.66"
Set-TextValue "E23" "  +0.96%  "
Set-TextValue "D24" "0.1385"
Set-TextValue "E24" "  -0.27%  "
Set-TextValue "D25" "8.404"
Set-TextValue "E25" "  +0.35%  "
Set-TextValue "D26" "17.69"
Set-TextValue "E26" "  +0.28%  "
Set-TextValue "D27" "1.374"
Set-TextValue "E27" "  +5.18%  "
Set-TextValue "D28" "1.460"
Set-TextValue "E28" "  -0.55%  "
Set-TextValue "D29" "0.05584"
Set-TextValue "E29" "  -0.83%  "
Set-TextValue "E30" "  -0.13%  "
Set-TextValue "D31" "4.070"
Set-TextValue "E31" "  +0.57%  "
Set-TextValue "D32" "1.835"
Set-TextValue "E32" "  -0.84%  "
Set-TextValue "D33" "1.163"
Set-TextValue "E33" "  +0.46%  "
Set-TextValue "D34" "0.6988"
Set-TextValue "E34" "  -1.53%  "
Set-TextValue "D35" "2.580"
Set-TextValue "E35" "  -0.23%  "
Set-TextValue "D36" "0.01804"
Set-TextValue "E36" "  +0.39%  "
Set-TextValue "D37" "1.231.84"
Set-TextValue "E37" "  +0.00%  "
Set-TextValue "D38" "2.714"
Set-TextValue "E38" "  -2.33%  "
Set-TextValue "D39" "6.374"
Set-TextValue "E39" "  -1.61%  "
Set-TextValue "D40" "0.9045"
Set-TextValue "E40" "  -0.46%  "
Set-TextValue "D41" "1.000"
Set-TextValue "E41" "  +0.01%  "
Set-TextValue "D42" "101.55"
Set-TextValue "E42" "  +0.20%  "
Set-TextValue "D43" "66.04"
Set-TextValue "B44" "BabyDogeCoin"
Set-TextValue "C44" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D44" "0.00000000121"
Set-TextValue "E44" "  +0.24%  "
Set-TextValue "B45" "Aptos"
Set-TextValue "C45" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D45" "7.194"
Set-TextValue "E45" "  +0.41%  "
Set-TextValue "B46" "TheSandbox"
Set-TextValue "C46" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D46" "0.4008"
Set-TextValue "E46" "  -0.09%  "
Set-TextValue "B47" "EnergySwap"
Set-TextValue "C47" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D47" "9.037"
Set-TextValue "E47" "  +0.22%  "
Set-TextValue "B48" "RenderToken"
Set-TextValue "C48" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D48" "1.678"
Set-TextValue "E48" "  -0.49%  "
Set-TextValue "B49" "Algorand"
Set-TextValue "C49" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D49" "0.1135"
Set-TextValue "E49" "  +1.08%  "
Set-TextValue "B50" "Cronos"
Set-TextValue "C50" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D50" "0.05703"
Set-TextValue "E50" "  -0.17%  "
Set-TextValue "B51" "Mantle"
Set-TextValue "C51" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D51" "0.4630"
Set-TextValue "E51" "  +0.11%  "
